$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_summaries")

$ws.Range("D2").Value = "{'batch_size': 16, 'epochs': 50, 'layers_struct': [{'units': 150, 'dropout': 0.4}, {'units': 100, 'dropout': 0.3}, {'units': 50, 'dropout': 0.2}]}"

$ws.Range("I2").Value = "rmse"
$ws.Range("J2").Value = 0.08771346807479929

$ws.Range("N2").Value = 1.390577358120786
$ws.Range("O2").Value = 0.1016039392577006
$ws.Range("P2").Value = 1.001877508455949
$ws.Range("Q2").Value = 0.5038096515647162
$ws.Range("R2").Value = 0.8316726753171663
$ws.Range("S2").Value = 0.06656464121876902
$ws.Range("T2").Value = 0.6604753503600755
$ws.Range("U2").Value = 0.8391876319974456
